$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force all touched cells to Text format first so numeric-looking strings
# (prices, percentages, hour values) keep their exact original text
# representation (leading/trailing zeros, "%" suffix, etc.) instead of being
# auto-converted into Excel numbers/doubles.
$cells = @(
    "D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5",
    "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8", "G8",
    "B9", "C9", "D9", "E9", "G9", "B10", "C10", "D10", "E10", "G10",
    "B11", "C11", "D11", "E11", "G11", "B12", "C12", "D12", "E12", "G12",
    "B13", "C13", "D13", "E13", "G13", "B14", "C14", "D14", "E14", "G14",
    "B15", "C15", "D15", "E15", "G15", "B16", "C16", "D16", "E16", "G16",
    "B17", "C17", "D17", "E17", "G17", "B18", "C18", "D18", "E18", "G18",
    "B19", "C19", "D19", "E19", "G19", "B20", "C20", "D20", "E20", "G20",
    "D21", "E21", "G21", "E22", "G22", "D23", "E23", "G23", "D24", "E24",
    "G24", "D25", "E25", "G25", "D26", "E26", "G26", "D27", "E27", "G27",
    "D28", "E28", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35",
    "G36", "G37", "G38", "G39", "D40", "E40", "G40", "D41", "E41", "G41",
    "E42", "G42", "D43", "E43", "G43", "D44", "E44", "G44", "D45", "E45",
    "G45", "E46", "G46", "E47", "G47", "D48", "E48", "G48", "D49", "E49",
    "G49", "D50", "E50", "G50"
)
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the new values cell by cell (grouped by row for readability).
# Row 2
$ws.Range("D2").Value = "245.24"
$ws.Range("E2").Value = "-0.67%"
$ws.Range("G2").Value = "20"

# Row 3
$ws.Range("D3").Value = "27.13"
$ws.Range("E3").Value = "2.70%"
$ws.Range("G3").Value = "20"

# Row 4
$ws.Range("D4").Value = "5.108"
$ws.Range("E4").Value = "0.46%"
$ws.Range("G4").Value = "20"

# Row 5
$ws.Range("D5").Value = "0.05706"
$ws.Range("G5").Value = "20"

# Row 6
$ws.Range("D6").Value = "6.506"
$ws.Range("E6").Value = "0.39%"
$ws.Range("G6").Value = "20"

# Row 7
$ws.Range("D7").Value = "0.8194"
$ws.Range("E7").Value = "0.75%"
$ws.Range("G7").Value = "20"

# Row 8
$ws.Range("D8").Value = "0.8637"
$ws.Range("E8").Value = "2.23%"
$ws.Range("G8").Value = "20"

# Row 9
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "0.1329"
$ws.Range("E9").Value = "-0.93%"
$ws.Range("G9").Value = "20"

# Row 10
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "0.06941"
$ws.Range("E10").Value = "-0.90%"
$ws.Range("G10").Value = "20"

# Row 11
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "0.02845"
$ws.Range("E11").Value = "-0.30%"
$ws.Range("G11").Value = "20"

# Row 12
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").Value = "0.09394"
$ws.Range("E12").Value = "-0.06%"
$ws.Range("G12").Value = "20"

# Row 13
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").Value = "0.001526"
$ws.Range("E13").Value = "1.11%"
$ws.Range("G13").Value = "20"

# Row 14
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").Value = "0.04029"
$ws.Range("E14").Value = "-13.28%"
$ws.Range("G14").Value = "20"

# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "0.0006018"
$ws.Range("E15").Value = "-93.95%"
$ws.Range("G15").Value = "20"

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.006203"
$ws.Range("E16").Value = "0.36%"
$ws.Range("G16").Value = "20"

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.511"
$ws.Range("E17").Value = "-2.70%"
$ws.Range("G17").Value = "20"

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "3.009"
$ws.Range("E18").Value = "-0.21%"
$ws.Range("G18").Value = "20"

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "2.318"
$ws.Range("E19").Value = "12.78%"
$ws.Range("G19").Value = "20"

# Row 20
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3165"
$ws.Range("E20").Value = "1.24%"
$ws.Range("G20").Value = "20"

# Row 21
$ws.Range("D21").Value = "0.03205"
$ws.Range("E21").Value = "0.23%"
$ws.Range("G21").Value = "20"

# Row 22
$ws.Range("E22").Value = "0.44%"
$ws.Range("G22").Value = "20"

# Row 23
$ws.Range("D23").Value = "3.583"
$ws.Range("E23").Value = "-4.30%"
$ws.Range("G23").Value = "20"

# Row 24
$ws.Range("D24").Value = "0.1374"
$ws.Range("E24").Value = "1.76%"
$ws.Range("G24").Value = "20"

# Row 25
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").Value = "-1.95%"
$ws.Range("G25").Value = "20"

# Row 26
$ws.Range("D26").Value = "0.004472"
$ws.Range("E26").Value = "-2.59%"
$ws.Range("G26").Value = "20"

# Row 27
$ws.Range("D27").Value = "0.00009899"
$ws.Range("E27").Value = "3.10%"
$ws.Range("G27").Value = "20"

# Row 28
$ws.Range("D28").Value = "0.0001937"
$ws.Range("E28").Value = "38.60%"
$ws.Range("G28").Value = "20"

# Row 29
$ws.Range("G29").Value = "20"

# Row 30
$ws.Range("G30").Value = "20"

# Row 31
$ws.Range("G31").Value = "20"

# Row 32
$ws.Range("G32").Value = "20"

# Row 33
$ws.Range("G33").Value = "20"

# Row 34
$ws.Range("G34").Value = "20"

# Row 35
$ws.Range("G35").Value = "20"

# Row 36
$ws.Range("G36").Value = "20"

# Row 37
$ws.Range("G37").Value = "20"

# Row 38
$ws.Range("G38").Value = "20"

# Row 39
$ws.Range("G39").Value = "20"

# Row 40
$ws.Range("D40").Value = "0.03728"
$ws.Range("E40").Value = "1.63%"
$ws.Range("G40").Value = "20"

# Row 41
$ws.Range("D41").Value = "0.005959"
$ws.Range("E41").Value = "-3.58%"
$ws.Range("G41").Value = "20"

# Row 42
$ws.Range("E42").Value = "-0.10%"
$ws.Range("G42").Value = "20"

# Row 43
$ws.Range("D43").Value = "0.002300"
$ws.Range("E43").Value = "-8.02%"
$ws.Range("G43").Value = "20"

# Row 44
$ws.Range("D44").Value = "0.009542"
$ws.Range("E44").Value = "6.67%"
$ws.Range("G44").Value = "20"

# Row 45
$ws.Range("D45").Value = "0.00005159"
$ws.Range("E45").Value = "-3.76%"
$ws.Range("G45").Value = "20"

# Row 46
$ws.Range("E46").Value = "-0.03%"
$ws.Range("G46").Value = "20"

# Row 47
$ws.Range("E47").Value = "-8.22%"
$ws.Range("G47").Value = "20"

# Row 48
$ws.Range("D48").Value = "0.002500"
$ws.Range("E48").Value = "-4.15%"
$ws.Range("G48").Value = "20"

# Row 49
$ws.Range("D49").Value = "0.00002100"
$ws.Range("E49").Value = "-0.03%"
$ws.Range("G49").Value = "20"

# Row 50
$ws.Range("D50").Value = "0.0002000"
$ws.Range("E50").Value = "-0.03%"
$ws.Range("G50").Value = "20"
